$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the cell contents that reference the workspace user path:
# user2 -> user6 (template path, project path, config path)
$ws.Range("B1").Value = "template /pub/home/user6/jmt_workspace/blocks/bl_1s24/bl_1s24.tsdl"
$ws.Range("J3").Value = "/pub/home/user6/jmt_workspace"
$ws.Range("K3").Value = "/pub/home/user6/jmt_workspace/workshop_config.sdl"

# Update the view: scroll so column G is the top-left visible column,
# and move the active selection to J4.
$ws.Application.ActiveWindow.ScrollColumn = 7
$ws.Range("J4").Select()
